$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before the existing row 454 (shifts old row 454.. down by 9)
$ws.Range("A454:A462").EntireRow.Insert()

$ws.Cells.Item(454,1).Value2 = 1574035200
$ws.Cells.Item(454,2).Value = "'2019-11-18"
$ws.Cells.Item(454,3).Value = "'0195"
$ws.Cells.Item(454,4).Value = "BINACOM"
$ws.Cells.Item(454,5).Value = 0.39
$ws.Cells.Item(454,6).Value = 0.39
$ws.Cells.Item(454,7).Value = 0.38
$ws.Cells.Item(454,8).Value = 0.38
$ws.Cells.Item(454,9).Value = 1212900

$ws.Cells.Item(455,1).Value2 = 1574121600
$ws.Cells.Item(455,2).Value = "'2019-11-19"
$ws.Cells.Item(455,3).Value = "'0195"
$ws.Cells.Item(455,4).Value = "BINACOM"
$ws.Cells.Item(455,5).Value = 0.385
$ws.Cells.Item(455,6).Value = 0.385
$ws.Cells.Item(455,7).Value = 0.375
$ws.Cells.Item(455,8).Value = 0.375
$ws.Cells.Item(455,9).Value = 1248600

$ws.Cells.Item(456,1).Value2 = 1574208000
$ws.Cells.Item(456,2).Value = "'2019-11-20"
$ws.Cells.Item(456,3).Value = "'0195"
$ws.Cells.Item(456,4).Value = "BINACOM"
$ws.Cells.Item(456,5).Value = 0.375
$ws.Cells.Item(456,6).Value = 0.395
$ws.Cells.Item(456,7).Value = 0.37
$ws.Cells.Item(456,8).Value = 0.38
$ws.Cells.Item(456,9).Value = 4850500

$ws.Cells.Item(457,1).Value2 = 1574294400
$ws.Cells.Item(457,2).Value = "'2019-11-21"
$ws.Cells.Item(457,3).Value = "'0195"
$ws.Cells.Item(457,4).Value = "BINACOM"
$ws.Cells.Item(457,5).Value = 0.375
$ws.Cells.Item(457,6).Value = 0.38
$ws.Cells.Item(457,7).Value = 0.375
$ws.Cells.Item(457,8).Value = 0.375
$ws.Cells.Item(457,9).Value = 796600

$ws.Cells.Item(458,1).Value2 = 1574380800
$ws.Cells.Item(458,2).Value = "'2019-11-22"
$ws.Cells.Item(458,3).Value = "'0195"
$ws.Cells.Item(458,4).Value = "BINACOM"
$ws.Cells.Item(458,5).Value = 0.375
$ws.Cells.Item(458,6).Value = 0.38
$ws.Cells.Item(458,7).Value = 0.335
$ws.Cells.Item(458,8).Value = 0.335
$ws.Cells.Item(458,9).Value = 3381100

$ws.Cells.Item(459,1).Value2 = 1574640000
$ws.Cells.Item(459,2).Value = "'2019-11-25"
$ws.Cells.Item(459,3).Value = "'0195"
$ws.Cells.Item(459,4).Value = "BINACOM"
$ws.Cells.Item(459,5).Value = 0.33
$ws.Cells.Item(459,6).Value = 0.35
$ws.Cells.Item(459,7).Value = 0.315
$ws.Cells.Item(459,8).Value = 0.32
$ws.Cells.Item(459,9).Value = 3973100

$ws.Cells.Item(460,1).Value2 = 1574726400
$ws.Cells.Item(460,2).Value = "'2019-11-26"
$ws.Cells.Item(460,3).Value = "'0195"
$ws.Cells.Item(460,4).Value = "BINACOM"
$ws.Cells.Item(460,5).Value = 0.325
$ws.Cells.Item(460,6).Value = 0.325
$ws.Cells.Item(460,7).Value = 0.305
$ws.Cells.Item(460,8).Value = 0.31
$ws.Cells.Item(460,9).Value = 2841100

$ws.Cells.Item(461,1).Value2 = 1574812800
$ws.Cells.Item(461,2).Value = "'2019-11-27"
$ws.Cells.Item(461,3).Value = "'0195"
$ws.Cells.Item(461,4).Value = "BINACOM"
$ws.Cells.Item(461,5).Value = 0.305
$ws.Cells.Item(461,6).Value = 0.31
$ws.Cells.Item(461,7).Value = 0.29
$ws.Cells.Item(461,8).Value = 0.29
$ws.Cells.Item(461,9).Value = 3975800

$ws.Cells.Item(462,1).Value2 = 1574899200
$ws.Cells.Item(462,2).Value = "'2019-11-28"
$ws.Cells.Item(462,3).Value = "'0195"
$ws.Cells.Item(462,4).Value = "BINACOM"
$ws.Cells.Item(462,5).Value = 0.29
$ws.Cells.Item(462,6).Value = 0.315
$ws.Cells.Item(462,7).Value = 0.285
$ws.Cells.Item(462,8).Value = 0.315
$ws.Cells.Item(462,9).Value = 2127500
